# Add 2022-Q3 data: insert a new per-quarter sheet right after "总计",
# and add a new summary row on "总计" for the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for "2022-Q3" and push the rest
#    down (the existing 2022-Q1 .. 2020-Q4 rows keep their values, the
#    row that used to be the last one now lands on row 8).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows("2:2").Insert()
$summary.Range("B2:D2").ClearFormats()

# A2 needs the same "index column" style ("s=2") the other A-cells use;
# copy it over from the row below (which keeps its original formatting).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Application.CutCopyMode = $false

$summaryData = @(
    @(0, "2022-Q3", 1, 0.64),
    @(1, "2022-Q1", 1, 0.02),
    @(2, "2021-Q4", 6, 1.8),
    @(3, "2021-Q3", 6, 1.24),
    @(4, "2021-Q2", 7, 2.32),
    @(5, "2021-Q1", 3, 0.65),
    @(6, "2020-Q4", 3, 1.43)
)

for ($i = 0; $i -lt $summaryData.Count; $i++) {
    $r = 2 + $i
    $row = $summaryData[$i]
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------
# 2. New "2022-Q3" worksheet: duplicate the existing "2022-Q1" sheet
#    (same header/layout/style) right after "总计", rename it, then
#    overwrite its single data row with the 2022-Q3 fund figures.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "100056"
$q3.Range("C2").Value = "富国低碳环保混合"

$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "21.80"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "83.39"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "2.93"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.6387"

$q3.Range("H2").Value = 10
